$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, A (WindowClassName), B (controlID), C (Module), D (Text)
# Ordered so that brand-new shared strings are first introduced in the same
# sequence as the target workbook, keeping sharedStrings.xml indices identical.
$entries = @(
    @(550, 'Edit', 22772, 'Beställning', 'Leverantörens namn'),
    @(555, 'Edit', 22797, 'Beställning', 'Projekt:'),
    @(556, 'Edit', 22795, 'Beställning', 'Beställningsdatum'),
    @(558, 'Edit', 22788, 'Beställning', 'Leverantörens ordernummer'),
    @(559, 'Button', 22792, 'Beställning', 'Skickad'),
    @(572, 'Edit', 22779, 'Beställning', 'Avvik. Namn'),
    @(573, 'Edit', 22780, 'Beställning', 'Avvik. Postadress'),
    @(574, 'Edit', 23172, 'Beställning', 'Avvik. Postadress 2'),
    @(575, 'Edit', 23731, 'Beställning', 'Avvik. GLN'),
    @(577, 'Edit', 22783, 'Beställning', 'Avvik. Postnummer'),
    @(578, 'Edit', 22782, 'Beställning', 'Avvik. Ort'),
    @(576, 'Edit', 22781, 'Beställning', 'Avvik. Besöksadress'),
    @(579, 'Edit', 23599, 'Beställning', 'Avvik. Landskod'),
    @(584, 'SafGrid', 21559, 'Beställning', 'Artiklar'),
    @(585, 'ComboBox', 22803, 'Beställning', 'Spårningsval'),
    @(549, 'Edit', 22770, 'Beställning', 'Leverantörsnummer'),
    @(551, 'Edit', 22793, 'Beställning', 'Betalningsvillkor'),
    @(552, 'Edit', 22807, 'Beställning', 'Leveransvillkor'),
    @(553, 'Edit', 22794, 'Beställning', 'Leveranssätt'),
    @(554, 'Edit', 22798, 'Beställning', 'Resultatenhet'),
    @(557, 'Edit', 22796, 'Beställning', 'Leveransdatum'),
    @(560, 'Edit', 22770, 'Beställning', 'Leverantörsnummer'),
    @(561, 'Edit', 22772, 'Beställning', 'Leverantörens namn'),
    @(562, 'Edit', 22774, 'Beställning', 'Postadress'),
    @(563, 'Edit', 23170, 'Beställning', 'Postadress 2'),
    @(564, 'Edit', 23730, 'Beställning', 'GLN'),
    @(565, 'Edit', 22776, 'Beställning', 'Postnummer'),
    @(566, 'Edit', 22777, 'Beställning', 'Ort'),
    @(567, 'Edit', 23598, 'Beställning', 'Landskod'),
    @(568, 'Edit', 22778, 'Beställning', 'Land'),
    @(569, 'Edit', 22804, 'Beställning', 'Vårt kundnummer'),
    @(570, 'Edit', 22786, 'Beställning', 'Språk'),
    @(571, 'Edit', 22787, 'Beställning', 'Valuta'),
    @(580, 'Edit', 22784, 'Beställning', 'Avvik. Landskod'),
    @(581, 'Edit', 22800, 'Beställning', 'Text 1'),
    @(582, 'Edit', 22801, 'Beställning', 'Text 2'),
    @(583, 'Edit', 22802, 'Beställning', 'Text 3'),
)

foreach ($e in $entries) {
    $r = $e[0]
    $ws.Cells.Item($r, 1).Value = $e[1]
    $ws.Cells.Item($r, 2).Value = $e[2]
    $ws.Cells.Item($r, 3).Value = $e[3]
    $ws.Cells.Item($r, 4).Value = $e[4]
}

$ws.Range("A1:E585").Select()